$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A82").Value = "legislative_turnover_heinsohn"
$ws.Range("B82").Value = "Legislative Turnover as measured by Heinsohn (2014)"

# Move selection to where the cursor lands after typing the new row
# (mirrors the author's final cursor position in the saved file).
$ws.Range("B83").Select()
